$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Source file (last N points)" todo (currently row 4) drops out of the
# still-open (Importance=1) group and lands at the bottom of it as the new
# row 10, with Importance flipped off (1 -> 0). Every row from 5-10 shifts
# up by one to fill the gap left behind.
#
# Helper: copy a row's A:D values (and B's cell style) from $src to $dst.
# Values are written first (plain assignment), then the source B cell's
# format is copied onto the destination B cell - PasteSpecial(xlPasteFormats)
# only reliably carries the style index when it runs after the value write.
function Copy-Row($src, $dst) {
    $ws.Range("A$dst").Value2 = $ws.Range("A$src").Value2
    $ws.Range("B$dst").Value2 = $ws.Range("B$src").Value2
    $ws.Range("C$dst").Value2 = $ws.Range("C$src").Value2
    $ws.Range("D$dst").Value2 = $ws.Range("D$src").Value2
    $ws.Range("B$src").Copy()
    $ws.Range("B$dst").PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

# Stash the original row 4 in an unused scratch row so it survives while
# rows 5-10 are shifted upward into 4-9.
$scratchRow = 25
Copy-Row 4 $scratchRow

# Shift rows 5-10 up into rows 4-9 (column E recalculates on its own since
# it's a formula).
for ($r = 5; $r -le 10; $r++) {
    $dst = $r - 1
    Copy-Row $r $dst
}

# Move the stashed original row 4 into row 10, with Importance (column B)
# flipped from 1 to 0 - this todo is no longer a current priority.
Copy-Row $scratchRow 10
$ws.Range("B10").Value2 = 0

# Clear the scratch row used as a temporary holding area.
$scratchRange = "A" + $scratchRow + ":D" + $scratchRow
$ws.Range($scratchRange).Clear()

# Selection moved from O17 to C9.
$ws.Range("C9").Select()
